# The workbook has two sheets: "Rent Data" and "test".
# The "test" sheet repeats a fixed 12-row x 7-column block of shared
# strings (source0..11, adress0..11, text0..11, price0..11, date0..11,
# url0..11, other0..11): rows 1-12, then again rows 13-24, then (after a
# blank gap at row 25) rows 26-37. This change appends that same 12-row
# block once more, leaving row 38 blank (matching the existing gap-row
# convention already used in this sheet), so the new copy lands on
# rows 39-50.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("test")

# Source block: the first occurrence of the repeating 12-row template.
$src = $ws.Range("A1:G12")

# Destination: rows 39-50 (row 38 is intentionally left empty, consistent
# with the existing blank-row gaps already present in this sheet).
$dst = $ws.Range("A39:G50")

$src.Copy($dst)
